# "Loan RBI, Variable Instalments"
#
# On the "Repayment schedule" sheet a new (blank) column is inserted
# before the old column N ("Late"), pushing the existing N/O/P columns
# (Late / heading / Outstanding) one place to the right (-> O/P/Q).
# The "Repayment schedule" tab also becomes the active tab/sheet
# (previously "Summary" was active), and the selection on that sheet
# moves to R7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make "Repayment schedule" the active sheet (this also flips
# tabSelected off of whichever sheet was active before - "Summary" -
# and on for this sheet, and updates the workbook's activeTab).
$ws.Activate()

# Insert a new blank column before column N; this shifts the existing
# N/O/P columns (and their column-width definitions) one to the right.
$ws.Columns("N:N").Insert()

# The newly inserted column takes on a plain width of 11 characters
# (no bestFit), matching column M's width but without the bestFit flag.
$ws.Columns("N:N").ColumnWidth = 10.17

# Update the active selection on the sheet.
$ws.Range("R7").Select() | Out-Null
